$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Forcing NumberFormat to Text before assignment prevents Excel from
    # auto-converting date-looking strings (e.g. "2023-06-30") into date
    # serial numbers. ClearFormats afterwards restores the cell's original
    # (default) style so only the value itself changes.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

# Row 2
$ws.Range("B2").Value = 89539

# Row 3
$ws.Range("A3").Value = 112183947
$ws.Range("B3").Value = 89704
$ws.Range("E3").Value = 1588
$ws.Range("F3").Value = "Violmussling"
$ws.Range("G3").Value = "Trichaptum laricinum"
$ws.Range("H3").Value = "(P.Karst.) Ryvarden"
$ws.Range("Q3").Value = 763391
$ws.Range("R3").Value = 7448820
Set-TextValue "Y3" "2023-06-30"
Set-TextValue "AA3" "2023-06-30"

# Row 4
$ws.Range("B4").Value = 89539

# Row 5
$ws.Range("B5").Value = 89557

# Row 6
$ws.Range("A6").Value = 112181997
$ws.Range("B6").Value = 89539
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = "Ullticka"
$ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 763401
$ws.Range("R6").Value = 7448827

# Row 7
$ws.Range("A7").Value = 112181983
$ws.Range("B7").Value = 89704
$ws.Range("Q7").Value = 763400
$ws.Range("R7").Value = 7448829
Set-TextValue "Y7" "2023-07-06"
Set-TextValue "AA7" "2023-07-06"
